$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-17 16:34:17"
$wsZhCn.Range("H4").Value = "2016-03-17 16:34:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-17 16:34:20"
$wsDeDe.Range("H4").Value = "2016-03-17 16:34:46"
